$wb = $excel.ActiveWorkbook

# --- Update the timestamps on the "data" sheet (column F) ---
$dataSheet = $wb.Worksheets.Item("data")

$dataSheet.Range("F2").Value = "2021-10-05 14:19:35.223925"
$dataSheet.Range("F3").Value = "2021-10-05 14:19:35.223932"
$dataSheet.Range("F4").Value = "2021-10-05 14:19:35.223935"
$dataSheet.Range("F5").Value = "2021-10-05 14:19:35.223937"
$dataSheet.Range("F6").Value = "2021-10-05 14:19:35.223940"
$dataSheet.Range("F7").Value = "2021-10-05 14:19:35.223943"

# --- Add the new "metadata" sheet (placed right after "data") ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"
$metaSheet.Outline.SummaryRow = 1
$metaSheet.Outline.SummaryColumn = 1

# Match the workbook's (openpyxl-default) page margins: 0.75/0.75/1/1/0.5/0.5in
# PageSetup margins are expressed in points (72pt = 1in).
$metaSheet.PageSetup.LeftMargin = 54
$metaSheet.PageSetup.RightMargin = 54
$metaSheet.PageSetup.TopMargin = 72
$metaSheet.PageSetup.BottomMargin = 72
$metaSheet.PageSetup.HeaderMargin = 36
$metaSheet.PageSetup.FooterMargin = 36

# Reuse the bold/bordered header style (cellXfs index "1") already present
# in the workbook by copy/pasting formats from the "data" sheet's header row
# and its styled A2 cell, instead of re-creating a brand-new style entry.
$dataSheet.Range("F1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)

$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Chondrodysplasia punctata"
$metaSheet.Range("C2").Value = 30

# "1.5" must be stored as literal text (not the number 1.5). A plain string
# assignment gets auto-coerced to a number, so force it via a leading
# quote-prefix, then paste-special the (default, un-styled) format from an
# always-empty cell to strip the resulting quotePrefix style back off.
$metaSheet.Range("D2").Value = "'1.5"
$dataSheet.Range("A1").Copy()
$metaSheet.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$metaSheet.Range("E2").Value = "2021-09-07T14:17:08.166044Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:19:35.220183"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/30/?format=json"

$dataSheet.Activate()
$dataSheet.Range("A1").Select()
